$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dual Boards")
$ws.Rows("33:33").Insert()
$ws.Range("E33").Value = "Lanyard, double ended"
$ws.Rows("33:33").Insert()
$ws.Range("E33").Value = "LED Diffuser"
$ws.Range("E33").NumberFormat = "General"
$ws.Range("E33").Style = "Currency"
